$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("B4").Value = 39647
$ws.Range("C4").Value = 3303.9
$ws.Range("D4").Value = 3684
$ws.Range("G4").Value = 2798
$ws.Range("H4").Value = 4689
$ws.Range("I4").Value = 12666
$ws.Range("K4").Value = 4955.9
$ws.Range("N4").Value = 0.3
$ws.Range("O4").Value = 70240
$ws.Range("P4").Value = 5853.3
$ws.Range("Q4").Value = 6683.4
$ws.Range("S4").Value = 63.8
$ws.Range("T4").Value = 3699
$ws.Range("U4").Value = 9791.799999999999
$ws.Range("V4").Value = 20340
$ws.Range("X4").Value = 7804.4
$ws.Range("AA4").Value = 0.7
$ws.Range("AB4").Value = 13864
$ws.Range("AC4").Value = 1155.3
$ws.Range("AD4").Value = 1062.8
$ws.Range("AF4").Value = 294
$ws.Range("AG4").Value = 1191
$ws.Range("AH4").Value = 1828.5
$ws.Range("AI4").Value = 3465
$ws.Range("AK4").Value = 1540.4
$ws.Range("AN4").Value = 0.6

# Row 5
$ws.Range("B5").Value = 898504
$ws.Range("C5").Value = 1804.2
$ws.Range("D5").Value = 2869.9
$ws.Range("F5").Value = 151.2
$ws.Range("G5").Value = 936.5
$ws.Range("H5").Value = 2169.2
$ws.Range("I5").Value = 24791
$ws.Range("K5").Value = 2240.7
$ws.Range("L5").Value = 401
$ws.Range("M5").Value = 80.5
$ws.Range("O5").Value = 393310
$ws.Range("P5").Value = 789.8
$ws.Range("Q5").Value = 2094
$ws.Range("S5").Value = 39.5
$ws.Range("T5").Value = 229
$ws.Range("U5").Value = 696.5
$ws.Range("V5").Value = 24554
$ws.Range("X5").Value = 961.6
$ws.Range("Y5").Value = 409
$ws.Range("Z5").Value = 82.09999999999999
$ws.Range("AB5").Value = 386451
$ws.Range("AC5").Value = 776
$ws.Range("AD5").Value = 768.7
$ws.Range("AF5").Value = 210
$ws.Range("AG5").Value = 633
$ws.Range("AH5").Value = 1110.2
$ws.Range("AI5").Value = 6193
$ws.Range("AK5").Value = 951.8
$ws.Range("AL5").Value = 406
$ws.Range("AM5").Value = 81.5
$ws.Range("AN5").Value = 1.2

# Row 6
$ws.Range("B6").Value = 793155
$ws.Range("C6").Value = 918
$ws.Range("D6").Value = 2530
$ws.Range("G6").Value = 212
$ws.Range("H6").Value = 1130
$ws.Range("I6").Value = 55081
$ws.Range("K6").Value = 1389.1
$ws.Range("L6").Value = 571
$ws.Range("M6").Value = 66.09999999999999
$ws.Range("N6").Value = 0.2
$ws.Range("O6").Value = 243540
$ws.Range("P6").Value = 281.9
$ws.Range("Q6").Value = 1851.2
$ws.Range("T6").Value = 32
$ws.Range("U6").Value = 184.2
$ws.Range("V6").Value = 48717
$ws.Range("X6").Value = 417
$ws.Range("Y6").Value = 584
$ws.Range("Z6").Value = 67.59999999999999
$ws.Range("AB6").Value = 464627
$ws.Range("AC6").Value = 537.8
$ws.Range("AD6").Value = 708.6
$ws.Range("AG6").Value = 300.5
$ws.Range("AH6").Value = 846
$ws.Range("AI6").Value = 8295
$ws.Range("AK6").Value = 750.6
$ws.Range("AL6").Value = 619
$ws.Range("AM6").Value = 71.59999999999999
$ws.Range("AN6").Value = 0.3

# Row 7
$ws.Range("B7").Value = 948515
$ws.Range("C7").Value = 519.7
$ws.Range("D7").Value = 1703.7
$ws.Range("G7").Value = 4
$ws.Range("H7").Value = 425
$ws.Range("I7").Value = 50844
$ws.Range("K7").Value = 1016.6
$ws.Range("L7").Value = 933
$ws.Range("M7").Value = 51.1
$ws.Range("O7").Value = 298239
$ws.Range("P7").Value = 163.4
$ws.Range("Q7").Value = 717
$ws.Range("U7").Value = 65
$ws.Range("V7").Value = 13956
$ws.Range("X7").Value = 313.9
$ws.Range("Y7").Value = 950
$ws.Range("Z7").Value = 52.1
$ws.Range("AA7").Value = -1
$ws.Range("AB7").Value = 538279
$ws.Range("AC7").Value = 294.9
$ws.Range("AD7").Value = 536.9
$ws.Range("AG7").Value = 35
$ws.Range("AH7").Value = 369
$ws.Range("AI7").Value = 5226
$ws.Range("AK7").Value = 523.6
$ws.Range("AL7").Value = 1028
$ws.Range("AM7").Value = 56.3
$ws.Range("AN7").Value = -1.1

# Row 8
$ws.Range("B8").Value = 143127
$ws.Range("C8").Value = 1022.3
$ws.Range("D8").Value = 4583.4
$ws.Range("G8").Value = 2
$ws.Range("H8").Value = 949.8
$ws.Range("I8").Value = 52263
$ws.Range("K8").Value = 1987.9
$ws.Range("L8").Value = 72
$ws.Range("M8").Value = 51.4
$ws.Range("N8").Value = -1
$ws.Range("O8").Value = 48925
$ws.Range("P8").Value = 349.5
$ws.Range("Q8").Value = 2125.4
$ws.Range("T8").Value = 1
$ws.Range("U8").Value = 112.5
$ws.Range("V8").Value = 24618
$ws.Range("X8").Value = 689.1
$ws.Range("Y8").Value = 71
$ws.Range("Z8").Value = 50.7
$ws.Range("AB8").Value = 61465
$ws.Range("AC8").Value = 439
$ws.Range("AD8").Value = 672.6
$ws.Range("AG8").Value = 53
$ws.Range("AH8").Value = 738.2
$ws.Range("AI8").Value = 3450
$ws.Range("AK8").Value = 768.3
$ws.Range("AL8").Value = 80
$ws.Range("AM8").Value = 57.1
$ws.Range("AN8").Value = -1
